$d = $word.ActiveDocument

# Locate the title-page "Date: January 25, 2024" paragraph (there is also an
# unrelated "Date: ____/____/____" signature-block placeholder near the end
# of the document that must NOT be touched) and the "Version 1" paragraph
# that immediately follows it.
$datePara = $null
$versionPara = $null

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd("`r")
    if ($t -eq "Date: January 25, 2024") {
        $datePara = $p
    }
    if ($t -eq "Version 1") {
        $versionPara = $p
    }
}

# Update "January 25, 2024" -> "February 19, 2024", scoped to the date
# paragraph so nothing else in the document can match.
if ($datePara -ne $null) {
    $rng = $datePara.Range
    $rng.Find.Execute("January 25, 2024", $true, $false, $false, $false, $false, $true, 1, $false, "February 19, 2024", 2) | Out-Null
} else {
    $d.Content.Find.Execute("Date: January 25, 2024", $true, $false, $false, $false, $false, $true, 1, $false, "Date: February 19, 2024", 2) | Out-Null
}

# Bump the version number "Version 1" -> "Version 2", scoped to its own
# paragraph.
if ($versionPara -ne $null) {
    $rng = $versionPara.Range
    $rng.Find.Execute("Version 1", $true, $false, $false, $false, $false, $true, 1, $false, "Version 2", 2) | Out-Null
} else {
    $d.Content.Find.Execute("Version 1", $true, $false, $false, $false, $false, $true, 1, $false, "Version 2", 2) | Out-Null
}
